$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-15 Sunday" "2023-10-16 Monday"

Replace-Text "24÷3=" "39÷8="
Replace-Text "35÷5=" "92÷5="
Replace-Text "63÷4=" "42÷5="
Replace-Text "84÷4=" "52÷8="
Replace-Text "33÷5=" "22÷3="

Replace-Text "84÷3=" "28÷5="
Replace-Text "18÷6=" "52÷3="
Replace-Text "28÷9=" "62÷3="
Replace-Text "58÷7=" "11÷3="
Replace-Text "61÷4=" "37÷2="

Replace-Text "76÷8=" "93÷3="
Replace-Text "44÷6=" "83÷2="
Replace-Text "89÷6=" "67÷7="
Replace-Text "40÷4=" "72÷2="
Replace-Text "15÷3=" "57÷5="

Replace-Text "23÷8=" "86÷3="
Replace-Text "73÷9=" "88÷9="
Replace-Text "82÷8=" "18÷8="
Replace-Text "38÷7=" "94÷6="
Replace-Text "42÷3=" "61÷7="

Replace-Text "23÷5=" "22÷8="
Replace-Text "38÷2=" "63÷7="
Replace-Text "93÷2=" "35÷6="
Replace-Text "11÷9=" "62÷4="
Replace-Text "57÷6=" "91÷4="
